# Update "想去人数" (number of people interested) counts for several
# events across the "展览", "演出", "本地生活" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4530
$ws1.Range("F5").Value = 3633
$ws1.Range("F6").Value = 1049
$ws1.Range("F9").Value = 358
$ws1.Range("F20").Value = 10387
$ws1.Range("F21").Value = 6045
$ws1.Range("F30").Value = 172
$ws1.Range("F32").Value = 3560
$ws1.Range("F35").Value = 477
$ws1.Range("F40").Value = 4846
$ws1.Range("F42").Value = 1132
$ws1.Range("F44").Value = 177
$ws1.Range("F45").Value = 99
$ws1.Range("F46").Value = 486

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 3560

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 439

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 439
$ws4.Range("F5").Value = 4530
$ws4.Range("F8").Value = 3633
$ws4.Range("F9").Value = 1049
$ws4.Range("F24").Value = 10387
$ws4.Range("F25").Value = 3560
$ws4.Range("F34").Value = 172
$ws4.Range("F36").Value = 3560
$ws4.Range("F42").Value = 4846
$ws4.Range("F44").Value = 1132
$ws4.Range("F46").Value = 99
$ws4.Range("F47").Value = 486
